# The workbook originally had:
#   A1 = 0            (numeric, styled: bold font + thin border + centered/top alignment)
#   A2 = long Python-ish single-line "questions = [...]" string (shared string)
#
# Target state:
#   A1 = the same question data, reformatted as pretty-printed JSON-ish text
#        (still stored as a plain string, no special styling)
#   row 2 is gone (used range shrinks back to just A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are analyzing the wages of employees in your company. Your colleague has already started the script. They have saved the information on administrative worker wages in the R vector called a and information on non-administrative worker wages in the R vector called b. Now you want to combine those two vectors into a single one. Which of the following operations will accomplish this?",
        "ques_type": 2,
        "options": [
            "a + b",
            "a | b",
            "c(a, b)",
            "a.extend(b)"
        ],
        "score": "c(a, b)"
    },
    {
        "title": "You are analyzing a car manufacturing company dataset stored in R. The dataset contains information about the characteristics of various cars. You want to build a classification model for predicting the car engine type given other characteristics. Before building the model, you need to identify what engine types exist in the dataset. Which of the following built-in functions in R should you use?",
        "ques_type": 2,
        "options": [
            "table",
            "select",
            "levels",
            "head"
        ],
        "score": "levels"
    },
    {
        "title": "You are analyzing the banking transaction data of your company which is stored in an R dataframe. You need to perform a generalized analysis but due to constraints, you have decided to sample the first 1,000 rows of the dataframe instead of analyzing all the data. Which method from the dplyr package can you use to return 1000 randomly sampled rows from the available dataframe? Note: Input a single word or expression only, e.g. method_name",
        "ques_type": null,
        "options": [],
        "score": null
    },
    {
        "title": "You want to create a chart containing boxplots of several variables available in your dataset using the ggplot function. You want to clearly define the orientation of the boxplots by setting the values to be represented on each axis. Which argument of the ggplot function should you use to fulfill the requirement?",
        "ques_type": 2,
        "options": [
            "label",
            "type",
            "aes",
            "coord"
        ],
        "score": "aes"
    }
]
'@

# Here-strings keep a single trailing newline before the closing '@ marker -
# strip it so the cell value ends exactly with the closing "]".
$text = $text.TrimEnd("`r", "`n")

# Row 2 (the old shared-string duplicate of the question text) is removed
# entirely, shrinking the sheet's used range back down to just row 1.
$ws.Range("A2").EntireRow.Delete()

# A1 becomes the reformatted text (replacing the old numeric 0), and loses
# the bold/bordered/centered styling it used to carry.
$ws.Range("A1").Value = $text
$ws.Range("A1").Style = "Normal"
